$wb = $excel.ActiveWorkbook

# NOTE: the Equipment sheet's C4/C5 multi-line descriptions are intentionally left
# untouched here. In the source edit, Excel re-serialized their embedded line
# breaks (LF) using its "_x000D_" escape notation, but that is a cosmetic
# re-encoding of the very same line-break content, not a content change -
# decoding "_x000D_" yields the same logical text/newlines already present.

# --- Summary Costs sheet: updated Use Cost ($/PMH) values in column I ---
$wsSummary = $wb.Worksheets.Item("Summary Costs")
$wsSummary.Range("I2").Value = 74.1269213331
$wsSummary.Range("I3").Value = 115.763080305
$wsSummary.Range("I4").Value = 74.55363122200001
$wsSummary.Range("I5").Value = 43.0491433378
$wsSummary.Range("I6").Value = 61.1466450328
$wsSummary.Range("I7").Value = 42.5960936825
$wsSummary.Range("I8").Value = 61.942622612
$wsSummary.Range("I9").Value = 101.813712625

# --- Operating Costs sheet: updated Tire/track ($/hr) values in column F ---
$wsOperating = $wb.Worksheets.Item("Operating Costs")
$wsOperating.Range("F2").Value = 3.83333333333
$wsOperating.Range("F3").Value = 2.14666666667
$wsOperating.Range("F4").Value = 0.383333333333
$wsOperating.Range("F5").Value = 0.383333333333
$wsOperating.Range("F6").Value = 4.44666666667
$wsOperating.Range("F7").Value = 0.383333333333
$wsOperating.Range("F8").Value = 5.75
$wsOperating.Range("F9").Value = 2.875
